$d = $word.ActiveDocument

# --- Step 1: remove everything from paragraph 19 ("tambien destaca...") through
# the end of the document (this covers "tambien destaca", "Internet movil y la
# adopcion del 5G", the duplicated "El acceso a Internet movil..." paragraph with
# the page break, "La penetracion...", "Telefonia, fija y movil",
# "Mientras el acceso...", "Por otro lado...", and the closing "El informe revela..."
# paragraph). Deleting via a range that runs to the document's content end avoids
# the special-cased "cannot delete the very last paragraph mark" behaviour.
$startPara = 19
$rng = $d.Range($d.Paragraphs($startPara).Range.Start, $d.Content.End)
$rng.Delete()

# --- Step 2: remove the "Te puede interesar..." paragraph (16) and the empty
# paragraph right before it (15), deleting one paragraph at a time (descending)
# so each deletion correctly collapses its own paragraph mark.
$d.Paragraphs(16).Range.Delete()
$d.Paragraphs(15).Range.Delete()

# --- Step 3: remove two of the three empty paragraphs that sit between
# "El crecimiento..." and "El numero de accesos fijos...", leaving just one.
$d.Paragraphs(12).Range.Delete()
$d.Paragraphs(11).Range.Delete()

# --- Step 4: remove the block of short paragraphs right after the title
# ("Hace 59 minsDigna Irene Urrea", "Internet movil en Colombia", "Facebook",
# "X", "LinkedIn", "WhatsApp"), one paragraph at a time (descending).
$d.Paragraphs(7).Range.Delete()
$d.Paragraphs(6).Range.Delete()
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(4).Range.Delete()
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(2).Range.Delete()

# --- Step 5: replace the remaining "El informe" paragraph text with the new
# mobile-internet paragraph text (assigning Range.Text keeps the paragraph's
# own trailing mark intact, so the paragraph count doesn't change). By this
# point paragraphs 2-7, 11-12, 15-16 and everything from 19 on have already
# been removed, so "El informe" (originally paragraph 17) is now paragraph 7.
$d.Paragraphs(7).Range.Text = "El acceso a Internet móvil también muestra un crecimiento notable, alcanzando los 48,1 millones de accesos, con un 85,4% operando en tecnología 4G y un 5,8% en 5G. La expansión de la red 5G es particularmente relevante, ya que en el último año se sumaron más de 640.000 nuevos accesos, lo que sugiere una transición progresiva hacia tecnologías más avanzadas. No obstante, la cobertura de 5G sigue siendo limitada a ciertos municipios, representando solo el 2% del total, lo que subraya la necesidad de mayores inversiones en infraestructura."

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
